# Auto-generated Excel COM-interop edit script
# Applies cell-value updates to match the target commit diff
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2568601.8
$ws.Range("I43").Value = 5131537.5
$ws.Range("J43").Value = 5666
$ws.Range("K43").Value = 5131537.5
$ws.Range("L43").Value = 5666
$ws.Range("M43").Value = -5131468.5
$ws.Range("N43").Value = -5804
$ws.Range("H100").Value = 767.875
$ws.Range("I100").Value = 800.5714
$ws.Range("K100").Value = 800.5714
$ws.Range("M100").Value = -259.5714
$ws.Range("H132").Value = 13496.471
$ws.Range("I132").Value = 2482.1035
$ws.Range("K132").Value = 7446.310500000001
$ws.Range("M132").Value = -4916.310500000001
$ws.Range("H133").Value = 94278
$ws.Range("J133").Value = 94278
$ws.Range("L133").Value = 94278
$ws.Range("N133").Value = -104398
$ws.Range("H141").Value = 9239.083000000001
$ws.Range("I141").Value = 9760.817999999999
$ws.Range("J141").Value = 3500
$ws.Range("K141").Value = 29282.454
$ws.Range("L141").Value = 10500
$ws.Range("M141").Value = -24102.454
$ws.Range("N141").Value = -20860

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4948.6665
$ws.Range("I45").Value = 4930.6665
$ws.Range("J45").Value = 4966.6665
$ws.Range("K45").Value = 4930.6665
$ws.Range("L45").Value = 4966.6665
$ws.Range("M45").Value = -4553.6665
$ws.Range("N45").Value = -5720.6665
$ws.Range("H63").Value = 2865
$ws.Range("I63").Value = 2865
$ws.Range("K63").Value = 2865
$ws.Range("M63").Value = -2179
$ws.Range("H66").Value = 2865
$ws.Range("I66").Value = 2865
$ws.Range("K66").Value = 14325
$ws.Range("M66").Value = -10893
$ws.Range("H74").Value = 8622361
$ws.Range("I74").Value = 11906000
$ws.Range("J74").Value = 2811.125
$ws.Range("K74").Value = 11906000
$ws.Range("L74").Value = 2811.125
$ws.Range("M74").Value = -11905126
$ws.Range("N74").Value = -4559.125
$ws.Range("H77").Value = 8622361
$ws.Range("I77").Value = 11906000
$ws.Range("J77").Value = 2811.125
$ws.Range("K77").Value = 59530000
$ws.Range("L77").Value = 14055.625
$ws.Range("M77").Value = -59525632
$ws.Range("N77").Value = -22791.625
$ws.Range("H97").Value = 4999
$ws.Range("I97").Value = 4999
$ws.Range("K97").Value = 4999
$ws.Range("M97").Value = -4503
$ws.Range("H102").Value = 320327.97
$ws.Range("I102").Value = 596836.5600000001
$ws.Range("J102").Value = 2343.05
$ws.Range("K102").Value = 596836.5600000001
$ws.Range("L102").Value = 2343.05
$ws.Range("M102").Value = -595214.5600000001
$ws.Range("N102").Value = -5587.05
$ws.Range("H110").Value = 4745.923
$ws.Range("I110").Value = 3078.2222
$ws.Range("J110").Value = 8498.25
$ws.Range("K110").Value = 3078.2222
$ws.Range("L110").Value = 8498.25
$ws.Range("M110").Value = -1033.2222
$ws.Range("N110").Value = -12588.25
$ws.Range("H122").Value = 3854.8635
$ws.Range("I122").Value = 2688.7273
$ws.Range("J122").Value = 5021
$ws.Range("K122").Value = 8066.1819
$ws.Range("L122").Value = 15063
$ws.Range("M122").Value = -5616.1819
$ws.Range("N122").Value = -19963
$ws.Range("H132").Value = 2782.1667
$ws.Range("I132").Value = 1509
$ws.Range("K132").Value = 4527
$ws.Range("M132").Value = -1997

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 806681.1
$ws.Range("I94").Value = 856973.7
$ws.Range("K94").Value = 856973.7
$ws.Range("M94").Value = -856522.7
$ws.Range("H134").Value = 8552.111000000001
$ws.Range("I134").Value = 7000
$ws.Range("K134").Value = 21000
$ws.Range("M134").Value = -18465

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 8000
$ws.Range("J25").Value = 7666.6665
$ws.Range("L25").Value = 7666.6665
$ws.Range("N25").Value = -8014.6665
$ws.Range("H94").Value = 3585.2
$ws.Range("I94").Value = 3198.2856
$ws.Range("J94").Value = 3923.75
$ws.Range("K94").Value = 3198.2856
$ws.Range("L94").Value = 3923.75
$ws.Range("M94").Value = -2747.2856
$ws.Range("N94").Value = -4825.75
$ws.Range("H100").Value = 51000.332
$ws.Range("J100").Value = 51000.332
$ws.Range("L100").Value = 51000.332
$ws.Range("N100").Value = -53164.332
$ws.Range("H134").Value = 3383.0833
$ws.Range("I134").Value = 3445.182
$ws.Range("K134").Value = 10335.546
$ws.Range("M134").Value = -7800.545999999998

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4991.615
$ws.Range("J3").Value = 12996
$ws.Range("L3").Value = 38988
$ws.Range("N3").Value = -39212
$ws.Range("H9").Value = 441
$ws.Range("I9").Value = 364.5
$ws.Range("J9").Value = 900
$ws.Range("K9").Value = 1093.5
$ws.Range("L9").Value = 2700
$ws.Range("M9").Value = -869.5
$ws.Range("N9").Value = -3148
$ws.Range("H37").Value = 166762480
$ws.Range("J37").Value = 166762480
$ws.Range("L37").Value = 500287440
$ws.Range("N37").Value = -500287664
$ws.Range("H131").Value = 8335906
$ws.Range("J131").Value = 6805387
$ws.Range("L131").Value = 20416161
$ws.Range("N131").Value = -20426241
$ws.Range("H133").Value = 18011.4
$ws.Range("I133").Value = 20797.777
$ws.Range("J133").Value = 13831.833
$ws.Range("K133").Value = 62393.33099999999
$ws.Range("L133").Value = 41495.499
$ws.Range("M133").Value = -57333.33099999999
$ws.Range("N133").Value = -51615.499

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2125
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2125
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2125
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -6465

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6308.3184
$ws.Range("I46").Value = 2150
$ws.Range("J46").Value = 6724.15
$ws.Range("K46").Value = 2150
$ws.Range("L46").Value = 6724.15
$ws.Range("M46").Value = -1962
$ws.Range("N46").Value = -7100.15
$ws.Range("H55").Value = 399.125
$ws.Range("I55").Value = 74.5
$ws.Range("K55").Value = 74.5
$ws.Range("M55").Value = 98.5
$ws.Range("H68").Value = 1265622
$ws.Range("I68").Value = 1895835.2
$ws.Range("J68").Value = 5195.3335
$ws.Range("K68").Value = 1895835.2
$ws.Range("L68").Value = 5195.3335
$ws.Range("M68").Value = -1895086.2
$ws.Range("N68").Value = -6693.3335
$ws.Range("H71").Value = 1265622
$ws.Range("I71").Value = 1895835.2
$ws.Range("J71").Value = 5195.3335
$ws.Range("K71").Value = 9479176
$ws.Range("L71").Value = 25976.6675
$ws.Range("M71").Value = -9475432
$ws.Range("N71").Value = -33464.6675
$ws.Range("H132").Value = 5597.579
$ws.Range("I132").Value = 5580.273
$ws.Range("K132").Value = 16740.819
$ws.Range("M132").Value = -14210.819

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 9250.25
$ws.Range("I2").Value = 10667
$ws.Range("K2").Value = 10667
$ws.Range("M2").Value = -10555
$ws.Range("H4").Value = 112.14286
$ws.Range("I4").Value = 116
$ws.Range("K4").Value = 116
$ws.Range("M4").Value = -3
$ws.Range("H132").Value = 37044670
$ws.Range("I132").Value = 6946382
$ws.Range("K132").Value = 20839146
$ws.Range("M132").Value = -20836616
